$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1904096666666667
$ws.Range("H2").Value = 0.571229
$ws.Range("I2").Value = 0.09975479936454949
$ws.Range("J2").Value = 0.1077412252037539
$ws.Range("M2").Value = 0.6598136666666666
$ws.Range("N2").Value = 1.979441
$ws.Range("O2").Value = 0.007704735356083927
$ws.Range("P2").Value = 0.008484678519943686
$ws.Range("Q2").Value = 0.1256349003321111
$ws.Range("R2").Value = 1.130714102989
$ws.Range("S2").Value = 0.0007685843296031029
$ws.Range("T2").Value = 0.0009141496591987061
$ws.Range("G3").Value = 0.1904096666666667
$ws.Range("H3").Value = 0.571229
$ws.Range("I3").Value = 0.09975479936454949
$ws.Range("J3").Value = 0.1077412252037539
$ws.Range("O3").Value = 0.7130079175842846
$ws.Range("P3").Value = 0.7851850431306702
$ws.Range("Q3").Value = 11.62644458527356
$ws.Range("R3").Value = 104.638001267462
$ws.Range("S3").Value = 0.07112596176395555
$ws.Range("T3").Value = 0.08459679855856077
$ws.Range("G4").Value = 0.1904096666666667
$ws.Range("H4").Value = 0.571229
$ws.Range("I4").Value = 0.09975479936454949
$ws.Range("J4").Value = 0.1077412252037539
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1563486666666667
$ws.Range("N4").Value = 0.469046
$ws.Range("O4").Value = 0.001825704984300993
$ws.Range("P4").Value = 0.002010519394650058
$ws.Range("Q4").Value = 0.02977029750377778
$ws.Range("R4").Value = 0.267932677534
$ws.Range("S4").Value = 0.0001821228344078035
$ws.Range("T4").Value = 0.0002166158228755069
$ws.Range("G5").Value = 0.1904096666666667
$ws.Range("H5").Value = 0.571229
$ws.Range("I5").Value = 0.09975479936454949
$ws.Range("J5").Value = 0.1077412252037539
$ws.Range("M5").Value = 23.6163295
$ws.Range("N5").Value = 47.232659
$ws.Range("O5").Value = 0.2757711427815902
$ws.Range("P5").Value = 0.2024581319964196
$ws.Range("Q5").Value = 4.496777427985166
$ws.Range("R5").Value = 26.980664567911
$ws.Range("S5").Value = 0.02750949501871006
$ws.Range("T5").Value = 0.02181308719375758
$ws.Range("G6").Value = 0.1904096666666667
$ws.Range("H6").Value = 0.571229
$ws.Range("I6").Value = 0.09975479936454949
$ws.Range("J6").Value = 0.1077412252037539
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.14477
$ws.Range("N6").Value = 0.43431
$ws.Range("O6").Value = 0.00169049929374041
$ws.Range("P6").Value = 0.001861626958316384
$ws.Range("Q6").Value = 0.02756560744333333
$ws.Range("R6").Value = 0.24809046699
$ws.Range("S6").Value = 0.0001686354178729872
$ws.Range("T6").Value = 0.000200573969361345
$ws.Range("I7").Value = 0.5298350724050007
$ws.Range("J7").Value = 0.5722539689365677
$ws.Range("M7").Value = 0.6598136666666666
$ws.Range("N7").Value = 1.979441
$ws.Range("O7").Value = 0.007704735356083927
$ws.Range("P7").Value = 0.008484678519943686
$ws.Range("Q7").Value = 0.6672939742056666
$ws.Range("R7").Value = 6.005645767851
$ws.Range("S7").Value = 0.004082239015252096
$ws.Range("T7").Value = 0.004855390958188617
$ws.Range("I8").Value = 0.5298350724050007
$ws.Range("J8").Value = 0.5722539689365677
$ws.Range("O8").Value = 0.7130079175842846
$ws.Range("P8").Value = 0.7851850431306702
$ws.Range("Q8").Value = 61.75239835969533
$ws.Range("S8").Value = 0.3777766016386082
$ws.Range("T8").Value = 0.4493252572811561
$ws.Range("I9").Value = 0.5298350724050007
$ws.Range("J9").Value = 0.5722539689365677
$ws.Range("K9").Value = 2
$ws.Range("L9").Value = 0.6666666666666666
$ws.Range("M9").Value = 0.1563486666666667
$ws.Range("N9").Value = 0.469046
$ws.Range("O9").Value = 0.001825704984300993
$ws.Range("P9").Value = 0.002010519394650058
$ws.Range("Q9").Value = 0.1581211915006666
$ws.Range("R9").Value = 1.423090723506
$ws.Range("S9").Value = 0.0009673225325472874
$ws.Range("T9").Value = 0.001150527703212441
$ws.Range("I10").Value = 0.5298350724050007
$ws.Range("J10").Value = 0.5722539689365677
$ws.Range("M10").Value = 23.6163295
$ws.Range("N10").Value = 47.232659
$ws.Range("O10").Value = 0.2757711427815902
$ws.Range("P10").Value = 0.2024581319964196
$ws.Range("Q10").Value = 23.8840678275415
$ws.Range("R10").Value = 143.304406965249
$ws.Range("S10").Value = 0.1461132234028936
$ws.Range("T10").Value = 0.1158574695784346
$ws.Range("I11").Value = 0.5298350724050007
$ws.Range("J11").Value = 0.5722539689365677
$ws.Range("K11").Value = 2
$ws.Range("L11").Value = 0.6666666666666666
$ws.Range("M11").Value = 0.14477
$ws.Range("N11").Value = 0.43431
$ws.Range("O11").Value = 0.00169049929374041
$ws.Range("P11").Value = 0.001861626958316384
$ws.Range("Q11").Value = 0.14641125749
$ws.Range("R11").Value = 1.31770131741
$ws.Range("S11").Value = 0.0008956858156995526
$ws.Range("T11").Value = 0.001065323415575861
$ws.Range("G12").Value = 0.08741566666666667
$ws.Range("H12").Value = 0.262247
$ws.Range("I12").Value = 0.04579668901430952
$ws.Range("J12").Value = 0.04946319792238989
$ws.Range("M12").Value = 0.6598136666666666
$ws.Range("N12").Value = 1.979441
$ws.Range("O12").Value = 0.007704735356083927
$ws.Range("P12").Value = 0.008484678519943686
$ws.Range("Q12").Value = 0.05767805154744444
$ws.Range("R12").Value = 0.519102463927
$ws.Range("S12").Value = 0.0003528513690401309
$ws.Range("T12").Value = 0.0004196793329398246
$ws.Range("G13").Value = 0.08741566666666667
$ws.Range("H13").Value = 0.262247
$ws.Range("I13").Value = 0.04579668901430952
$ws.Range("J13").Value = 0.04946319792238989
$ws.Range("O13").Value = 0.7130079175842846
$ws.Range("P13").Value = 0.7851850431306702
$ws.Range("Q13").Value = 5.337614534896222
$ws.Range("R13").Value = 48.038530814066
$ws.Range("S13").Value = 0.03265340186634792
$ws.Range("T13").Value = 0.03883776319407258
$ws.Range("G14").Value = 0.08741566666666667
$ws.Range("H14").Value = 0.262247
$ws.Range("I14").Value = 0.04579668901430952
$ws.Range("J14").Value = 0.04946319792238989
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.1563486666666667
$ws.Range("N14").Value = 0.469046
$ws.Range("O14").Value = 0.001825704984300993
$ws.Range("P14").Value = 0.002010519394650058
$ws.Range("Q14").Value = 0.01366732292911111
$ws.Range("R14").Value = 0.123005906362
$ws.Range("S14").Value = 0.000083611243397907416023337457
$ws.Range("T14").Value = 0.000099446718744379344217566852
$ws.Range("G15").Value = 0.08741566666666667
$ws.Range("H15").Value = 0.262247
$ws.Range("I15").Value = 0.04579668901430952
$ws.Range("J15").Value = 0.04946319792238989
$ws.Range("M15").Value = 23.6163295
$ws.Range("N15").Value = 47.232659
$ws.Range("O15").Value = 0.2757711427815902
$ws.Range("P15").Value = 0.2024581319964196
$ws.Range("Q15").Value = 2.064437187462167
$ws.Range("R15").Value = 12.386623124773
$ws.Range("S15").Value = 0.01262940526508923
$ws.Range("T15").Value = 0.01001422665393624
$ws.Range("G16").Value = 0.08741566666666667
$ws.Range("H16").Value = 0.262247
$ws.Range("I16").Value = 0.04579668901430952
$ws.Range("J16").Value = 0.04946319792238989
$ws.Range("K16").Value = 2
$ws.Range("L16").Value = 0.6666666666666666
$ws.Range("M16").Value = 0.14477
$ws.Range("N16").Value = 0.43431
$ws.Range("O16").Value = 0.00169049929374041
$ws.Range("P16").Value = 0.001861626958316384
$ws.Range("Q16").Value = 0.01265516606333333
$ws.Range("R16").Value = 0.11389649457
$ws.Range("S16").Value = 0.000077419270434339430841941265
$ws.Range("T16").Value = 0.000092082022696859993022081015
$ws.Range("G17").Value = 0.42447
$ws.Range("H17").Value = 0.84894
$ws.Range("I17").Value = 0.2223779938672774
$ws.Range("J17").Value = 0.1601211348241683
$ws.Range("M17").Value = 0.6598136666666666
$ws.Range("N17").Value = 1.979441
$ws.Range("O17").Value = 0.007704735356083927
$ws.Range("P17").Value = 0.008484678519943686
$ws.Range("Q17").Value = 0.28007110709
$ws.Range("R17").Value = 1.68042664254
$ws.Range("S17").Value = 0.001713363591764226
$ws.Range("T17").Value = 0.001358576353231628
$ws.Range("G18").Value = 0.42447
$ws.Range("H18").Value = 0.84894
$ws.Range("I18").Value = 0.2223779938672774
$ws.Range("J18").Value = 0.1601211348241683
$ws.Range("O18").Value = 0.7130079175842846
$ws.Range("P18").Value = 0.7851850431306702
$ws.Range("Q18").Value = 25.91820583222
$ws.Range("R18").Value = 155.50923499332
$ws.Range("S18").Value = 0.1585572703238782
$ws.Range("T18").Value = 0.1257247201530465
$ws.Range("G19").Value = 0.42447
$ws.Range("H19").Value = 0.84894
$ws.Range("I19").Value = 0.2223779938672774
$ws.Range("J19").Value = 0.1601211348241683
$ws.Range("K19").Value = 2
$ws.Range("L19").Value = 0.6666666666666666
$ws.Range("M19").Value = 0.1563486666666667
$ws.Range("N19").Value = 0.469046
$ws.Range("O19").Value = 0.001825704984300993
$ws.Range("P19").Value = 0.002010519394650058
$ws.Range("Q19").Value = 0.06636531854
$ws.Range("R19").Value = 0.39819191124
$ws.Range("S19").Value = 0.0004059966118023439
$ws.Range("T19").Value = 0.0003219266470573673
$ws.Range("G20").Value = 0.42447
$ws.Range("H20").Value = 0.84894
$ws.Range("I20").Value = 0.2223779938672774
$ws.Range("J20").Value = 0.1601211348241683
$ws.Range("M20").Value = 23.6163295
$ws.Range("N20").Value = 47.232659
$ws.Range("O20").Value = 0.2757711427815902
$ws.Range("P20").Value = 0.2024581319964196
$ws.Range("Q20").Value = 10.024423382865
$ws.Range("R20").Value = 40.09769353146
$ws.Range("S20").Value = 0.06132543349825652
$ws.Range("T20").Value = 0.03241782584964797
$ws.Range("G21").Value = 0.42447
$ws.Range("H21").Value = 0.84894
$ws.Range("I21").Value = 0.2223779938672774
$ws.Range("J21").Value = 0.1601211348241683
$ws.Range("K21").Value = 2
$ws.Range("L21").Value = 0.6666666666666666
$ws.Range("M21").Value = 0.14477
$ws.Range("N21").Value = 0.43431
$ws.Range("O21").Value = 0.00169049929374041
$ws.Range("P21").Value = 0.001861626958316384
$ws.Range("Q21").Value = 0.06145052190000001
$ws.Range("R21").Value = 0.3687031314000001
$ws.Range("S21").Value = 0.0003759298415760415
$ws.Range("T21").Value = 0.0002980858211848842
$ws.Range("G22").Value = 0.1951446666666667
$ws.Range("H22").Value = 0.585434
$ws.Range("I22").Value = 0.102235445348863
$ws.Range("J22").Value = 0.1104204731131201
$ws.Range("M22").Value = 0.6598136666666666
$ws.Range("N22").Value = 1.979441
$ws.Range("O22").Value = 0.007704735356083927
$ws.Range("P22").Value = 0.008484678519943686
$ws.Range("Q22").Value = 0.1287591180437778
$ws.Range("R22").Value = 1.158832062394
$ws.Range("S22").Value = 0.0007876970504243709
$ws.Range("T22").Value = 0.0009368822163849092
$ws.Range("G23").Value = 0.1951446666666667
$ws.Range("H23").Value = 0.585434
$ws.Range("I23").Value = 0.102235445348863
$ws.Range("J23").Value = 0.1104204731131201
$ws.Range("O23").Value = 0.7130079175842846
$ws.Range("P23").Value = 0.7851850431306702
$ws.Range("Q23").Value = 11.91556443971689
$ws.Range("R23").Value = 107.240079957452
$ws.Range("S23").Value = 0.07289468199149475
$ws.Range("T23").Value = 0.08670050394383419
$ws.Range("G24").Value = 0.1951446666666667
$ws.Range("H24").Value = 0.585434
$ws.Range("I24").Value = 0.102235445348863
$ws.Range("J24").Value = 0.1104204731131201
$ws.Range("K24").Value = 2
$ws.Range("L24").Value = 0.6666666666666666
$ws.Range("M24").Value = 0.1563486666666667
$ws.Range("N24").Value = 0.469046
$ws.Range("O24").Value = 0.001825704984300993
$ws.Range("P24").Value = 0.002010519394650058
$ws.Range("Q24").Value = 0.03051060844044444
$ws.Range("R24").Value = 0.274595475964
$ws.Range("S24").Value = 0.000186651762145651
$ws.Range("T24").Value = 0.0002220025027603632
$ws.Range("G25").Value = 0.1951446666666667
$ws.Range("H25").Value = 0.585434
$ws.Range("I25").Value = 0.102235445348863
$ws.Range("J25").Value = 0.1104204731131201
$ws.Range("M25").Value = 23.6163295
$ws.Range("N25").Value = 47.232659
$ws.Range("O25").Value = 0.2757711427815902
$ws.Range("P25").Value = 0.2024581319964196
$ws.Range("Q25").Value = 4.608600748167667
$ws.Range("R25").Value = 27.651604489006
$ws.Range("S25").Value = 0.02819358559664076
$ws.Range("T25").Value = 0.02235552272064317
$ws.Range("G26").Value = 0.1951446666666667
$ws.Range("H26").Value = 0.585434
$ws.Range("I26").Value = 0.102235445348863
$ws.Range("J26").Value = 0.1104204731131201
$ws.Range("K26").Value = 2
$ws.Range("L26").Value = 0.6666666666666666
$ws.Range("M26").Value = 0.14477
$ws.Range("N26").Value = 0.43431
$ws.Range("O26").Value = 0.00169049929374041
$ws.Range("P26").Value = 0.001861626958316384
$ws.Range("Q26").Value = 0.02825109339333334
$ws.Range("R26").Value = 0.25425984054
$ws.Range("S26").Value = 0.0001728289481574892
$ws.Range("T26").Value = 0.0002055617294974338
